$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Places"
$ws.Range("B1").Value = "Coordinates"
$ws.Range("A4").Value = "PreToria"
$ws.Range("B4").Value = "25.7479° S, 28.2293° E"
$ws.Range("C1").Value = "TimesTamp"
$ws.Range("A2").Value = "astana"
$ws.Range("A5").Value = "almaty"
$ws.Range("B5").Value = "43.2220° N, 76.8512° E"
$ws.Range("A3").Value = "murmask"
$ws.Range("A6").Value = "tromso"
$ws.Range("B2").Value = "51.1605° N, 71.4704° E"
$ws.Range("A7").Value = "naples"
$ws.Range("B6").Value = "69.6492° N, 18.9553° E"
$ws.Range("A8").Value = "milan"
$ws.Range("A9").Value = "belgrade"
$ws.Range("B8").Value = "45.4642° N, 9.1900° E"
$ws.Range("A10").Value = "kingston"
$ws.Range("B9").Value = "44.7866° N, 20.4489° E"
$ws.Range("B10").Value = "44.2312° N, 76.4860° W"

$ws.Columns.Item(3).ColumnWidth = 19.86

$ws.Range("A10").Select()
